$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Headers contain the "BTec_Logo-Orange" picture; its OOXML name attribute
# (wp:docPr/name and pic:cNvPr/name) changes from image2.jpg -> image1.jpg
for ($i = 1; $i -le 2; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Range.InlineShapes.Count -gt 0) {
        $shp = $hdr.Range.InlineShapes.Item(1)
        if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
            $shp.Name = "image1.jpg"
        }
    }
}

# Footers contain the Pearson Edexcel logo picture; its OOXML name attribute
# changes from image1.png -> image2.png
for ($i = 1; $i -le 2; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Range.InlineShapes.Count -gt 0) {
        $shp = $ftr.Range.InlineShapes.Item(1)
        if ($shp.AlternativeText -like "*PearsonLogo.png") {
            $shp.Name = "image2.png"
        }
    }
}
